$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in building depth (voxel_depth, column D) for rows 3 through 20 with
# value 13, matching the already-populated D2 cell, and copy D2's style so
# the formatting (bold font / fill / border) matches the rest of the table.
$ws.Range("D2").Copy() | Out-Null
$ws.Range("D3:D20").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

for ($r = 3; $r -le 20; $r++) {
    $ws.Cells.Item($r, 4).Value = 13
}

# Update the last selected cell to match the saved workbook view.
$ws.Range("P16").Select() | Out-Null
